$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.645200000000001
$ws.Range("B8").Value = 5.602799999999998
$ws.Range("B10").Value = 5.828499999999998
$ws.Range("D11").Value = -7.459999999999996
$ws.Range("B12").Value = 5.126599999999999
$ws.Range("D12").Value = -6.050099999999998
$ws.Range("D15").Value = -8.7158
$ws.Range("D17").Value = -8.441299999999993
$ws.Range("B18").Value = 6.517699999999995
$ws.Range("B25").Value = 6.419499999999996
$ws.Range("D26").Value = -7.083800000000001
$ws.Range("D27").Value = -8.492700000000001
$ws.Range("D28").Value = -8.519599999999999
$ws.Range("D32").Value = -6.371999999999997
$ws.Range("B37").Value = 8.621300000000002
$ws.Range("D37").Value = -7.620399999999998
$ws.Range("D41").Value = -7.931599999999996
$ws.Range("D47").Value = -7.769499999999997
$ws.Range("D51").Value = -8.681899999999999
$ws.Range("B55").Value = 5.532299999999997
$ws.Range("D65").Value = -7.955900000000002
$ws.Range("B68").Value = 4.694799999999996
$ws.Range("D73").Value = -8.754899999999997
$ws.Range("B77").Value = 8.520300000000002
$ws.Range("B78").Value = 8.740999999999998
$ws.Range("B79").Value = 8.794900000000005
$ws.Range("B80").Value = 8.781699999999999
$ws.Range("B81").Value = 5.800000000000004
$ws.Range("B82").Value = 6.186299999999997
$ws.Range("B84").Value = 6.580700000000004
$ws.Range("D84").Value = -7.303899999999999
$ws.Range("D85").Value = -8.495799999999999
$ws.Range("D89").Value = -8.542799999999996
$ws.Range("D93").Value = -6.611799999999993
$ws.Range("D95").Value = -7.583800000000002
$ws.Range("D98").Value = -6.822099999999997
$ws.Range("D99").Value = -8.467700000000004
$ws.Range("B101").Value = 5.820099999999996
$ws.Range("D101").Value = -7.664099999999996
$ws.Range("B102").Value = 7.912500000000005
$ws.Range("D102").Value = -8.008500000000003
